$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Some Price values are numeric-looking strings (e.g. "1.00", "8.03") that must
# stay as text (matching the original inlineStr cells), so we force the cell
# to Text format before assigning those specific values.

$ws.Range("D2").Value = "69.225.13"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").Value = "3.734.76"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.52"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.37"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").Value = "3.734.13"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.42"
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.15"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000247"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "4.356.73"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "3.737.05"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").Value = "69.204.22"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.29"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.04"
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.02"
$ws.Range("E21").Value = "  +19.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.22"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.726"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("E24").Value = "  +10.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.81"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.30"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.19"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +2.68%  "
$ws.Range("E31").Value = "  +6.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("E32").Value = "  +4.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.73"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").Value = "3.880.67"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "3.667.92"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.87"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("E42").Value = "  +8.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "432.90"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.69"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.47"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.41"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.86"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "2.782.40"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("E51").Value = "  +0.67%  "
